$d = $word.ActiveDocument
$tbl = $d.Tables.Item(1)

$tbl.Cell(1,1).Range.Text = "79÷3="
$tbl.Cell(1,2).Range.Text = "90÷5="
$tbl.Cell(1,3).Range.Text = "58÷7="
$tbl.Cell(1,4).Range.Text = "23÷8="
$tbl.Cell(1,5).Range.Text = "33÷3="
$tbl.Cell(5,1).Range.Text = "72÷9="
$tbl.Cell(5,2).Range.Text = "16÷8="
$tbl.Cell(5,3).Range.Text = "99÷8="
$tbl.Cell(5,4).Range.Text = "61÷7="
$tbl.Cell(5,5).Range.Text = "71÷6="
$tbl.Cell(9,1).Range.Text = "96÷7="
$tbl.Cell(9,2).Range.Text = "67÷4="
$tbl.Cell(9,3).Range.Text = "49÷5="
$tbl.Cell(9,4).Range.Text = "11÷3="
$tbl.Cell(9,5).Range.Text = "77÷7="
$tbl.Cell(13,1).Range.Text = "19÷3="
$tbl.Cell(13,2).Range.Text = "35÷8="
$tbl.Cell(13,3).Range.Text = "87÷7="
$tbl.Cell(13,4).Range.Text = "19÷7="
$tbl.Cell(13,5).Range.Text = "61÷2="
$tbl.Cell(17,1).Range.Text = "83÷3="
$tbl.Cell(17,2).Range.Text = "54÷2="
$tbl.Cell(17,3).Range.Text = "84÷6="
$tbl.Cell(17,4).Range.Text = "19÷5="
$tbl.Cell(17,5).Range.Text = "84÷4="
